$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$m = $ws.Range("Y1").MergeArea.Address()
Write-Output $m
